# "se terminaron las mediciones" — finish the measurement table.
#
# The Y (Frecuencia), Z (Vinput), AA (Voutput) and AC (phase) columns for
# the second measurement block (circuito 1 - caso 1x1, low-pass sweep)
# only had readings starting at 1000 Hz. The missing low-frequency
# readings (100 Hz .. 950 Hz) are filled in, which pushes the existing
# 1000 Hz..3000 Hz readings further down the table (rows 16-30). The very
# first row of the block (row 5, the 1000 Hz reading before the edit) no
# longer holds a measurement, so its inputs are cleared out (its gain
# formula is left in place and now evaluates to #DIV/0!, same as the
# still-unmeasured rows below the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 no longer has a measurement - clear the raw inputs, keep the
# formula in AB5 (it will naturally recompute to #DIV/0!).
$ws.Range("Y5").Clear()
$ws.Range("Z5").Clear()
$ws.Range("AA5").Clear()
$ws.Range("AC5").Clear()

# New / shifted measurements, 100 Hz through 3000 Hz.
$measurements = @{
  6  = @(100,  0.199, 0.02,  180)
  7  = @(200,  0.201, 0.021, 170)
  8  = @(300,  0.201, 0.021, 165)
  9  = @(400,  0.201, 0.021, 160)
  10 = @(500,  0.201, 0.021, 157)
  11 = @(600,  0.201, 0.021, 151)
  12 = @(700,  0.201, 0.021, 152)
  13 = @(800,  0.201, 0.021, 140)
  14 = @(900,  0.202, 0.021, 135)
  15 = @(950,  0.202, 0.021, 130)
  16 = @(1000, 0.481, 0.053, 120)
  17 = @(1100, 0.48,  0.052, 111)
  18 = @(1200, 0.48,  0.05,  102)
  19 = @(1300, 0.479, 0.048, 93)
  20 = @(1400, 0.477, 0.045, 85)
  21 = @(1500, 0.477, 0.042, 78)
  22 = @(1600, 0.475, 0.039, 70)
  23 = @(1700, 0.474, 0.036, 65)
  24 = @(1750, 0.473, 0.035, 62)
  25 = @(1780, 0.473, 0.034, 61)
  26 = @(1800, 0.473, 0.034, 59)
  27 = @(1900, 0.472, 0.031, 54)
  28 = @(2000, 0.469, 0.03,  50)
  29 = @(2500, 0.46,  0.021, 34)
  30 = @(3000, 0.451, 0.016, 20)
}

foreach ($row in ($measurements.Keys | Sort-Object { [int]$_ })) {
  $vals = $measurements[$row]
  $ws.Range("Y$row").Value = $vals[0]
  $ws.Range("Z$row").Value = $vals[1]
  $ws.Range("AA$row").Value = $vals[2]
  $ws.Range("AC$row").Value = $vals[3]
}

# Row 38 (AB5:AB38's old tail) never had a measurement and is now fully
# removed from the sheet instead of being left as a dangling #DIV/0!.
$ws.Range("AB38").Clear()

$wb.Save()
